$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix malformed path in E4: backslash after 02_PRODUCTOS should be a forward slash
$ws.Range("E4").Value = "s3://foa-prod-comp-fenomenologico-bucket/foa_puj_curada/P0016/02_PRODUCTOS/PRODUCTO 1/Anexos producto 1/Mapa 1-3. PoligonosThiessen.pdf"

# Update the active cell selection to match the saved view state
$ws.Range("E9").Select()
